# The "Recorded By" column (G) lists reviewers as a comma-separated
# string. Every cell whose value is literally "dnasr281@gmail.com, System"
# needs to have the two names swapped to "System, dnasr281@gmail.com"
# (cells that already read "System, dnasr281@gmail.com" are left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed++
    }
}

Write-Output "Replaced $changed cell(s) in column G (Recorded By)"
